$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "66.574.02"
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.581.13"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  -0.03%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "608.26"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.46%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "146.32"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +1.32%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "3.581.11"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +0.64%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("E10").Value = "  -0.10%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "8.00"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -0.47%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.415"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +1.02%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "4.194.86"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("E14").Value = "  +0.37%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "30.04"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -0.71%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "3.565.11"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +0.23%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "66.685.54"
$cell.Style = "Normal"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "11.37"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -2.37%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "6.28"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.92%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "14.97"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.11%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "432.55"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.98%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.621"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +2.00%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "79.12"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.54%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "3.730.34"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  +0.36%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "9.29"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.55%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "8.06"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  +1.05%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +0.02%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.582.55"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +0.73%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "25.47"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("E34").Value = "  -3.02%  "
$ws.Range("E35").Value = "  -1.77%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "7.84"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +0.02%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "1.71"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -2.39%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "5.63"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -0.53%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "174.37"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +1.51%  "
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("E42").Value = "  -1.38%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.893"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("E44").Value = "  +0.75%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "45.68"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("E46").Value = "  -0.07%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.55"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +6.22%  "
$ws.Range("E48").Value = "  -2.04%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "24.90"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -4.24%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "23.76"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +3.39%  "
$ws.Range("E51").Value = "  +0.78%  "
